$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.867.44'
$ws.Range("E2").Value = '  +0.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.792.99'
$ws.Range("E3").Value = '  -0.42%  '

# Row 4
$ws.Range("E4").Value = '  -0.40%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.61'
$ws.Range("E5").Value = '  -1.12%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.24'
$ws.Range("E6").Value = '  -2.19%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.790.48'
$ws.Range("E7").Value = '  -0.51%  '

# Row 8
$ws.Range("E8").Value = '  -0.14%  '

# Row 9
$ws.Range("E9").Value = '  -1.39%  '

# Row 10
$ws.Range("E10").Value = '  -2.61%  '

# Row 11
$ws.Range("E11").Value = '  -0.95%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.65'
$ws.Range("E12").Value = '  +5.57%  '

# Row 13
$ws.Range("E13").Value = '  -3.36%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.17'
$ws.Range("E14").Value = '  -2.47%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.426.10'
$ws.Range("E15").Value = '  -0.42%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.789.53'
$ws.Range("E16").Value = '  -0.89%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.884.66'
$ws.Range("E17").Value = '  +0.06%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.15'
$ws.Range("E18").Value = '  -1.92%  '

# Row 19
$ws.Range("E19").Value = '  +2.09%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.02'
$ws.Range("E20").Value = '  -1.03%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '458.28'
$ws.Range("E21").Value = '  -0.88%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.45'
$ws.Range("E22").Value = '  -4.36%  '

# Row 23
$ws.Range("E23").Value = '  -1.24%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.99'
$ws.Range("E24").Value = '  -0.44%  '

# Row 25
$ws.Range("E25").Value = '  -5.32%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.87'
$ws.Range("E26").Value = '  -1.88%  '

# Row 27
$ws.Range("E27").Value = '  -1.53%  '

# Row 28
$ws.Range("E28").Value = '  -0.12%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.90'
$ws.Range("E29").Value = '  -1.19%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.936.91'
$ws.Range("E30").Value = '  -0.53%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.22'
$ws.Range("E31").Value = '  -2.46%  '

# Row 32
$ws.Range("E32").Value = '  -1.93%  '

# Row 33
$ws.Range("E33").Value = '  -7.71%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.92'
$ws.Range("E34").Value = '  -2.35%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.03%  '

# Row 36
$ws.Range("E36").Value = '  -1.56%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0991'
$ws.Range("E37").Value = '  -0.96%  '

# Row 38
$ws.Range("E38").Value = '  +4.85%  '

# Row 39
$ws.Range("E39").Value = '  +0.22%  '

# Row 40
$ws.Range("E40").Value = '  -1.94%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.18'
$ws.Range("E41").Value = '  -6.07%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("E43").Value = '  +0.09%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.77'
$ws.Range("E44").Value = '  +1.51%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.13'
$ws.Range("E45").Value = '  -2.21%  '

# Row 46
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.89'
$ws.Range("E46").Value = '  +2.08%  '

# Row 47
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.294'
$ws.Range("E47").Value = '  -2.28%  '

# Row 48
$ws.Range("E48").Value = '  -0.89%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  -0.88%  '

# Row 50
$ws.Range("E50").Value = '  -0.77%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.44'
$ws.Range("E51").Value = '  -5.80%  '
